$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Geyser" data row with the new "TabelFan" asset type
# (assetTypeName / assetTypeCode / assetTypeDescription)
$ws.Range("A2").Value = "TabelFan"
$ws.Range("B2").Value = "TFN"
$ws.Range("C2").Value = "Table fans "

# Leave the sheet with a fresh selection, as the user did after editing
$ws.Range("A3:XFD15").Select()

# Reflect the resized application window in the saved view state
$excel.ActiveWindow.Width = 23370
$excel.ActiveWindow.Height = 10740
